$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Fgf1", "Fgfr1", "ECs", 2, 0.6666666666666666, 1.161357666666667, 3.484073, 0.1270850363824361, 0.1270850363824361, 3, 1, 3.76917, 11.30751, 0.02686645020528053, 0.02686645020528053, 4.377354476470001, 39.39619028823, 0.003414323801804983, 0.003414323801804984),
    @("ECs", "Fgf1", "Fgfr1", "FAPs", 2, 0.6666666666666666, 1.161357666666667, 3.484073, 0.1270850363824361, 0.1270850363824361, 3, 1, 82.48060333333333, 247.44181, 0.5879175050094569, 0.5879175050094571, 95.78948103245889, 862.10532929213, 0.07471551751399789, 0.07471551751399791),
    @("ECs", "Fgf1", "Fgfr1", "M2", 2, 0.6666666666666666, 1.161357666666667, 3.484073, 0.1270850363824361, 0.1270850363824361, 3, 1, 0.5716056666666667, 1.714817, 0.004074375838860061, 0.004074375838860062, 0.6638386232934445, 5.974547609641, 0.0005177922017172494, 0.0005177922017172496),
    @("ECs", "Fgf1", "Fgfr1", "sCs", 2, 0.6666666666666666, 1.161357666666667, 3.484073, 0.1270850363824361, 0.1270850363824361, 3, 1, 53.47143866666666, 160.414316, 0.3811416689464024, 0.3811416689464024, 62.09946524322978, 558.895187189068, 0.04843740286491596, 0.04843740286491597),
    @("FAPs", "Fgf1", "Fgfr1", "ECs", 3, 1, 5.168173666666667, 15.504521, 0.565542861868062, 0.565542861868062, 3, 1, 3.76917, 11.30751, 0.02686645020528053, 0.02686645020528053, 19.47972513919, 175.31752625271, 0.01519412913733013, 0.01519412913733013),
    @("FAPs", "Fgf1", "Fgfr1", "FAPs", 3, 1, 5.168173666666667, 15.504521, 0.565542861868062, 0.565542861868062, 3, 1, 82.48060333333333, 247.44181, 0.5879175050094569, 0.5879175050094571, 426.2740821581123, 3836.46673942301, 0.3324925483253789, 0.3324925483253791),
    @("FAPs", "Fgf1", "Fgfr1", "M2", 3, 1, 5.168173666666667, 15.504521, 0.565542861868062, 0.565542861868062, 3, 1, 0.5716056666666667, 1.714817, 0.004074375838860061, 0.004074375838860062, 2.954157354184111, 26.587416187657, 0.002304234172235005, 0.002304234172235005),
    @("FAPs", "Fgf1", "Fgfr1", "sCs", 3, 1, 5.168173666666667, 15.504521, 0.565542861868062, 0.565542861868062, 3, 1, 53.47143866666666, 160.414316, 0.3811416689464024, 0.3811416689464024, 276.3496812358484, 2487.147131122636, 0.2155519502331179, 0.2155519502331179),
    @("sCs", "Fgf1", "Fgfr1", "ECs", 3, 1, 2.808898333333333, 8.426695, 0.3073721017495019, 0.3073721017495019, 3, 1, 3.76917, 11.30751, 0.02686645020528053, 0.02686645020528053, 10.58721533105, 95.28493797945001, 0.008257997266145413, 0.008257997266145413),
    @("sCs", "Fgf1", "Fgfr1", "FAPs", 3, 1, 2.808898333333333, 8.426695, 0.3073721017495019, 0.3073721017495019, 3, 1, 82.48060333333333, 247.44181, 0.5879175050094569, 0.5879175050094571, 231.6796292353278, 2085.11666311795, 0.1807094391700801, 0.1807094391700801),
    @("sCs", "Fgf1", "Fgfr1", "M2", 3, 1, 2.808898333333333, 8.426695, 0.3073721017495019, 0.3073721017495019, 3, 1, 0.5716056666666667, 1.714817, 0.004074375838860061, 0.004074375838860062, 1.605582204423889, 14.450239839815, 0.001252349464907807, 0.001252349464907807),
    @("sCs", "Fgf1", "Fgfr1", "sCs", 3, 1, 2.808898333333333, 8.426695, 0.3073721017495019, 0.3073721017495019, 3, 1, 53.47143866666666, 160.414316, 0.3811416689464024, 0.3811416689464024, 150.1958349517355, 1351.76251456562, 0.1171523158483686, 0.1171523158483686)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowVals[$j]
    }
}
